$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove two workers from the EC database (YUDELKA GOMEZ VILLEGAS - old row 16,
# and ALFONSO LUIS BELLO ARROYO - old row 18). Delete from the bottom row up so
# row numbers of not-yet-deleted rows remain stable.
$ws.Rows("18:18").Delete()
$ws.Rows("16:16").Delete()

# The remaining last data row (now row 18, previously row 20 - DARWIS DAVID
# MARTINEZ BARBOZA) gets a new "Periodo Mora" / "Valor Mora" entry as part 1 of
# the new account statement additions.
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940

# Update the summary figures: total overdue value and worker count.
$ws.Range("E11").Value = 118720
$ws.Range("C13").Value = 2
